# Apply the changes described by the diff to the workbook.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Reporte de Formatos"
$ws2 = $wb.Worksheets.Item(2)   # "Tabla_393859"

# -----------------------------------------------------------------
# Sheet1 ("Reporte de Formatos") - row 8 content changes
# -----------------------------------------------------------------

# I8: used to contain a long paragraph (now emptied), keep a plain
# bordered style with no special alignment.
$ws1.Range("I8").ClearContents()
$ws1.Range("I8").Style = "Normal"
$ws1.Range("I8").Borders.LineStyle = 1

# D8: gains a value of 1 (style already correct: bordered + left align)
$ws1.Range("D8").Value = 1

# -----------------------------------------------------------------
# Sheet1 row heights
# -----------------------------------------------------------------
$ws1.Rows.Item(3).RowHeight = 96.75
$ws1.Rows.Item(5).RowHeight = 6
$ws1.Rows.Item(8).RowHeight = 35.25

# -----------------------------------------------------------------
# Sheet1 column widths
# -----------------------------------------------------------------
$ws1.Columns.Item(5).ColumnWidth = 110.193
$ws1.Columns.Item(9).ColumnWidth = 35.586

# -----------------------------------------------------------------
# Sheet2 ("Tabla_393859") - new data rows (4-8)
# -----------------------------------------------------------------
$data = @(
    @(1, 1000, "Servicios Personales", 75416997, -2042402.12, 73374594.88, 73374594.88, 73374594.88, 0),
    @(1, 2000, "Materiales y Suministros", 1715034, -977449.3, 737584.7, 737584.7, 737584.7, 0),
    @(1, 3000, "Servicios Generales", 9157271, -1890478.13, 7266792.87, 7266792.87, 7258793.13, 0),
    @(1, 4000, "Transferencia, Asignaciones, Subsidios y Otras Ayudas", 84000, -72806, 11194, 11194, 11194, 0),
    @(1, 5000, "Bienes Muebles, Inmuebles e Intangibles", 736000, -537232, 198768, 198768, 198768, 0)
)

$r = 4
foreach ($row in $data) {
    $rng = $ws2.Range("A" + $r + ":I" + $r)
    $rng.Borders.LineStyle = 1
    $rng.HorizontalAlignment = -4131

    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
    $ws2.Cells.Item($r, 6).Value = $row[5]
    $ws2.Cells.Item($r, 7).Value = $row[6]
    $ws2.Cells.Item($r, 8).Value = $row[7]
    $ws2.Cells.Item($r, 9).Value = $row[8]

    $r = $r + 1
}

# -----------------------------------------------------------------
# Sheet2 column widths
# -----------------------------------------------------------------
$ws2.Columns.Item(3).ColumnWidth = 57.92547
$ws2.Columns.Item(6).ColumnWidth = 18.92547
$ws2.Columns.Item(7).ColumnWidth = 19.59
$ws2.Columns.Item(8).ColumnWidth = 20.92547

# -----------------------------------------------------------------
# E8: gains a hyperlink to the budget PDF report
# -----------------------------------------------------------------
$ws1.Range("E8").Style = "Normal"
$ws1.Range("E8").Borders.LineStyle = 1
$ws1.Hyperlinks.Add($ws1.Range("E8"), "https://www.upp.edu.mx/leygralcontabilidad/mc/02-edospres/05-informacion-presupuestaria/2021/a_septiembre_2021/06.estadoanaliticopresupuestoegresos-cap-gto_ex.pdf")
$ws1.Range("E8").WrapText = $true

# -----------------------------------------------------------------
# Selections / view state
# -----------------------------------------------------------------
$ws2.Range("B12").Select()
$ws1.Select()
$ws1.Range("B10").Select()
